# MAI_holdings.xlsx - "Add files via upload" update
#
#  - bump the "as of" date in the confidential disclosure note (A10) from
#    2021-05-25 to 2021-05-26
#  - refresh the Weight (D) / Percent Change (E) figures for rows 2-7
#
# The sheet ships protected, so it has to be unprotected for the writes to
# go through, then protection is restored once the edits are in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Disclosure text: bump the "as of" date --------------------------------
$cell = $ws.Range("A10")
$cell.Value2 = $cell.Value2.Replace("2021-05-25", "2021-05-26")
# Re-run the automatic row height now that the text changed, instead of
# leaving an explicit (and incorrect) manual row height behind.
$ws.Rows.Item(10).AutoFit()

# --- Weight (D) / Percent Change (E) refresh --------------------------------
$ws.Range("D2").Value = 0.4782899344827649
$ws.Range("E2").Value = -0.0003901677721418473

$ws.Range("D3").Value = 0.3390297457946013
$ws.Range("E3").Value = 0.001788908765653074

$ws.Range("D4").Value = 0.09747353804732489
$ws.Range("E4").Value = 0.004408393581379011

$ws.Range("D5").Value = 0.05354083424747259
$ws.Range("E5").Value = 0.0003439578078421324

$ws.Range("D6").Value = 0.0316659474278362
$ws.Range("E6").Value = 0.006978772899098651

$ws.Range("E7").Value = 0.001088986929087543

$ws.Protect()
